$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 9.237899999999994
$ws.Range("B6").Value = 6.097800000000001
$ws.Range("B7").Value = 5.137999999999997
$ws.Range("C7").Value = -14.15959999999999
$ws.Range("C12").Value = -11.0212
$ws.Range("C15").Value = -14.65129999999998
$ws.Range("B16").Value = 6.792199999999997
$ws.Range("B20").Value = 9.522999999999996
$ws.Range("C20").Value = -11.80630000000001
$ws.Range("C21").Value = -11.8738
$ws.Range("C22").Value = -12.6137
$ws.Range("C23").Value = -12.1373
$ws.Range("B28").Value = 5.525799999999998
$ws.Range("B29").Value = 4.833300000000002
$ws.Range("C29").Value = -10.46360000000001
$ws.Range("B32").Value = 7.152699999999997
$ws.Range("C34").Value = -11.32600000000001
$ws.Range("B40").Value = 9.203299999999986
$ws.Range("C42").Value = -12.41900000000001
$ws.Range("C43").Value = -13.15459999999998
$ws.Range("C44").Value = -13.79829999999999
$ws.Range("C45").Value = -13.97829999999998
$ws.Range("B46").Value = 5.743099999999997
$ws.Range("C46").Value = -14.09719999999999
$ws.Range("C50").Value = -14.10609999999999
$ws.Range("B51").Value = 5.795700000000002
$ws.Range("C51").Value = -12.1999
$ws.Range("B52").Value = 5.209000000000001
$ws.Range("B57").Value = 5.240899999999995
$ws.Range("B59").Value = 5.4477
$ws.Range("B62").Value = 5.713000000000002
$ws.Range("B66").Value = 5.677200000000002
$ws.Range("C66").Value = -11.4373
$ws.Range("C67").Value = -11.5975
$ws.Range("B73").Value = 8.748899999999999
$ws.Range("B74").Value = 9.196099999999989
$ws.Range("C79").Value = -11.8657
$ws.Range("C84").Value = -13.63979999999999
$ws.Range("B92").Value = 5.569399999999995
$ws.Range("C92").Value = -11.5104
$ws.Range("C97").Value = -12.20420000000001
$ws.Range("B100").Value = 5.845300000000003

$wb.Save()